$d = $word.ActiveDocument

# --- Step 1: merge the heading paragraph's two runs into one --------------
# In the original document the heading paragraph reads:
#   "Adding a style here, with a unique word that’s easy to search for"
#   <bookmarkStart _GoBack/><bookmarkEnd/>
#   ": helicopter."
# as two separate <w:r> runs split by the (hidden) "_GoBack" bookmark that
# Word drops at the last edit location. Re-finding/replacing the whole
# sentence with itself collapses it back down to a single run and drops
# that now-stale bookmark (it is about to be relocated in step 2, just like
# Word relocates "_GoBack" to the newest edit point).
$curly_apos = [char]0x2019
$headingText = "Adding a style here, with a unique word that" + $curly_apos + "s easy to search for: helicopter."
$headingRange = $d.Content
$headingRange.Find.Execute($headingText, $false, $false, $false, $false, $false, $true, 1, $false, $headingText, 2)

# --- Step 2: add a blank paragraph, then a new paragraph with the typed ---
# --- text, carrying the relocated "_GoBack" bookmark ----------------------
$endRange = $d.Content
$endRange.Collapse(0)
$wNs = 'xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"'
$blankParagraphXml = "<w:p $wNs/>"
$typedParagraphXml = "<w:p $wNs><w:r><w:t>sdsadasd</w:t></w:r><w:bookmarkStart w:id=" + '"0"' + " w:name=" + '"_GoBack"' + "/><w:bookmarkEnd w:id=" + '"0"' + "/></w:p>"
$newParagraphsXml = $blankParagraphXml + $typedParagraphXml
$endRange.InsertXML($newParagraphsXml)
